$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.621.78'
$ws.Range("E2").Value = '  +3.85%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.31'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.51'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5288'
$ws.Range("E7").Value = '  -1.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3778'
$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07530'
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.118'
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("E14").Value = '  +6.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.189'
$ws.Range("E15").Value = '  +0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.795.66'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.27'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001067'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06469'
$ws.Range("E19").Value = '  -0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.29'
$ws.Range("E21").Value = '  +1.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.932'
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.627.42'
$ws.Range("E23").Value = '  +3.74%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.85'
$ws.Range("E26").Value = '  +3.41%  '

$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.377'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.002.91'
$ws.Range("E29").Value = '  +0.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.89'
$ws.Range("E30").Value = '  +1.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.117'
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.699'
$ws.Range("E33").Value = '  +0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.680'
$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2274'
$ws.Range("E35").Value = '  +8.99%  '

$ws.Range("E36").Value = '  +8.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.920'
$ws.Range("E37").Value = '  +2.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02315'
$ws.Range("E38").Value = '  +1.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.064'
$ws.Range("E39").Value = '  +1.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.46'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6290'
$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.208'
$ws.Range("E42").Value = '  +5.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.45'
$ws.Range("E45").Value = '  +1.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5916'
$ws.Range("E46").Value = '  +0.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.663'
$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.84'
$ws.Range("E48").Value = '  +4.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.975'
$ws.Range("E49").Value = '  +3.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.160'
$ws.Range("E50").Value = '  +2.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06926'
$ws.Range("E51").Value = '  +2.67%  '
